# BAU Fraction of CCS Potential Achieved by Industry - update
$wb = $excel.ActiveWorkbook

# --- About sheet: update date in C1, make it the active/selected sheet ---
$about = $wb.Worksheets.Item("About")
$about.Range("C1").Value = "4/5/2024"

# --- BAU Emissions sheet: replace the " : NoSettings" suffix with " : test"
#     in the long list of series labels in column A, update one row of
#     numeric data, and update the viewport/selection ---
$bau = $wb.Worksheets.Item("BAU Emissions")
$bau.Range("A1:A300").Replace(": NoSettings", ": test") | Out-Null

$bau.Range("M94").Value  = 1001080
$bau.Range("N94").Value  = 2002150
$bau.Range("O94").Value  = 3003230
$bau.Range("P94").Value  = 4004300
$bau.Range("Q94").Value  = 5005380
$bau.Range("R94").Value  = 5005380
$bau.Range("S94").Value  = 5005380
$bau.Range("T94").Value  = 5005380
$bau.Range("U94").Value  = 5005380
$bau.Range("V94").Value  = 5005380
$bau.Range("W94").Value  = 5005380
$bau.Range("X94").Value  = 5005380
$bau.Range("Y94").Value  = 5005380
$bau.Range("Z94").Value  = 5005380
$bau.Range("AA94").Value = 5005380
$bau.Range("AB94").Value = 5005380
$bau.Range("AC94").Value = 5005380
$bau.Range("AD94").Value = 5005380
$bau.Range("AE94").Value = 5005380

# Activate "About" last so it ends up the selected/active tab in the saved file
$about.Activate()
$about.Range("E29").Select()

# Restore the BAU Emissions sheet's on-screen selection (A30:AE280)
$bau.Activate()
$bau.Range("A30:AE280").Select()

# Leave the final active sheet as "About" (matches target workbook view)
$about.Activate()
